{"js": "// The document contains a placeholder date \"Fecha: dd/mm/aaaa\" that the\n// author filled in with the real inspection date (09/06/2021).\nconst body = context.document.body;\n\nconst results = body.search(\"Fecha: dd/mm/aaaa\", { matchCase: false, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"Fecha: 09/06/2021\", Word.InsertLocation.replace);\n  }\n} else {\n  // Fallback: the placeholder text might already have been partially\n  // edited; try to target just the \"dd/mm/aaaa\" token instead.\n  const fallback = body.search(\"dd/mm/aaaa\", { matchCase: false });\n  fallback.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < fallback.items.length; i++) {\n    fallback.items[i].insertText(\"09/06/2021\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a placeholder date \"Fecha: dd/mm/aaaa\" that the\n# author filled in with the real inspection date (09/06/2021).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$found = $find.Execute(\n    [ref]\"Fecha: dd/mm/aaaa\",  # FindText\n    [ref]$false,               # MatchCase\n    [ref]$false,               # MatchWholeWord\n    [ref]$false,               # MatchWildcards\n    [ref]$false,               # MatchSoundsLike\n    [ref]$false,               # MatchAllWordForms\n    [ref]$true,                # Forward\n    [ref]1,                    # Wrap (wdFindContinue)\n    [ref]$false,               # Format\n    [ref]\"Fecha: 09/06/2021\",  # ReplaceWith\n    [ref]2                     # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    # Fallback: the placeholder text might already have been partially\n    # edited; try to target just the \"dd/mm/aaaa\" token instead.\n    $find2 = $d.Content.Find\n    $find2.Execute(\n        [ref]\"dd/mm/aaaa\",\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$true,\n        [ref]1,\n        [ref]$false,\n        [ref]\"09/06/2021\",\n        [ref]2\n    )\n}\n"}
